# Obligation Suite task: swap the highlighted ("current") row in the
# Test Suite sheet from row 19 (WOR Suite) to row 12 (Obligation Suite),
# and update the scrolled/selected view position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the yellow highlight formatting between row 12 and row 19 ---
# Row 19 currently carries the yellow fill; row 3 (and row 12) carry the
# plain/no-fill format. Use Copy + PasteSpecial(Formats) so we reuse the
# workbook's existing cell styles instead of inventing new ones.
$ws.Range("A19:C19").Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null

$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A19:C19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Swap the Runmode ("Y"/"N") values between the two rows ---
$ws.Range("C12").Value = "Y"
$ws.Range("C19").Value = "N"

# --- Move the frozen-pane scroll position / active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B10").Select() | Out-Null
